$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_1_8_0"
$ws.Range("B2").Value = 0.9835144447384616
$ws.Range("C2").Value = 0.9831426913416406
$ws.Range("D2").Value = 0.9964837803326917
$ws.Range("E2").Value = 0.9906393015909399
$ws.Range("F2").Value = 1.98816287279599
$ws.Range("G2").Value = 0.9449729881512847
$ws.Range("H2").Value = 0.221428842062661
$ws.Range("I2").Value = 0.6044811801729448

$ws.Range("A3").Value = "model_1_8_1"
$ws.Range("B3").Value = 0.9861920806992862
$ws.Range("C3").Value = 0.9837610579032268
$ws.Range("D3").Value = 0.993607526522157
$ws.Range("E3").Value = 0.9896035472651509
$ws.Range("F3").Value = 1.665239178706346
$ws.Range("G3").Value = 0.9103091097518566
$ws.Range("H3").Value = 0.4025567609655702
$ws.Range("I3").Value = 0.6713665737474522

$ws.Range("A4").Value = "model_1_8_2"
$ws.Range("B4").Value = 0.988246992434232
$ws.Range("C4").Value = 0.9841905484273437
$ws.Range("D4").Value = 0.9898364786591878
$ws.Range("E4").Value = 0.9880703637333854
$ws.Range("F4").Value = 1.417416211661757
$ws.Range("G4").Value = 0.8862330871682592
$ws.Range("H4").Value = 0.6400330396587453
$ws.Range("I4").Value = 0.7703742065332971

$ws.Range("A5").Value = "model_1_8_3"
$ws.Range("B5").Value = 0.9898082646033566
$ws.Range("C5").Value = 0.9844689645616889
$ws.Range("D5").Value = 0.985475974479419
$ws.Range("E5").Value = 0.9861972358566303
$ws.Range("F5").Value = 1.229126323226811
$ws.Range("G5").Value = 0.8706258670743616
$ws.Range("H5").Value = 0.9146294763696339
$ws.Range("I5").Value = 0.891334256742784

$ws.Range("A6").Value = "model_1_8_4"
$ws.Range("B6").Value = 0.9909794727868443
$ws.Range("C6").Value = 0.9846285357356932
$ws.Range("D6").Value = 0.9807632736803805
$ws.Range("E6").Value = 0.9841078875093551
$ws.Range("F6").Value = 1.087878267593672
$ws.Range("G6").Value = 0.8616807589211182
$ws.Range("H6").Value = 1.211404985198328
$ws.Range("I6").Value = 1.026257069075994

$ws.Range("A7").Value = "model_1_8_5"
$ws.Range("B7").Value = 0.991843800343518
$ws.Range("C7").Value = 0.9846962938565862
$ws.Range("D7").Value = 0.9758801786853707
$ws.Range("E7").Value = 0.9818981376607989
$ws.Range("F7").Value = 0.9836401069220362
$ws.Range("G7").Value = 0.857882430536118
$ws.Range("H7").Value = 1.518910821787514
$ws.Range("I7").Value = 1.168954989462934

$ws.Range("A8").Value = "model_1_8_24"
$ws.Range("B8").Value = 0.992414361998941
$ws.Range("C8").Value = 0.9826677854025788
$ws.Range("D8").Value = 0.9169414997716347
$ws.Range("E8").Value = 0.9539185408834218
$ws.Range("F8").Value = 0.9148302014043631
$ws.Range("G8").Value = 0.9715948702927987
$ws.Range("H8").Value = 5.230488783172969
$ws.Range("I8").Value = 2.975779538407078

$ws.Range("A9").Value = "model_1_8_6"
$ws.Range("B9").Value = 0.9924678621361126
$ws.Range("C9").Value = 0.9846946919713304
$ws.Range("D9").Value = 0.9709637392753373
$ws.Range("E9").Value = 0.9796412262047719
$ws.Range("F9").Value = 0.9083780689328351
$ws.Range("G9").Value = 0.8579722276874597
$ws.Range("H9").Value = 1.828516474630103
$ws.Range("I9").Value = 1.314698441593021

$ws.Range("A10").Value = "model_1_8_23"
$ws.Range("B10").Value = 0.99249752464824
$ws.Range("C10").Value = 0.9827417899625712
$ws.Range("D10").Value = 0.9182812519616631
$ws.Range("E10").Value = 0.9545673880053641
$ws.Range("F10").Value = 0.9048007611388326
$ws.Range("G10").Value = 0.967446384220072
$ws.Range("H10").Value = 5.146119828967108
$ws.Range("I10").Value = 2.933879259508661

$ws.Range("A11").Value = "model_1_8_22"
$ws.Range("B11").Value = 0.9925860354408741
$ws.Range("C11").Value = 0.982822479958409
$ws.Range("D11").Value = 0.9197562937200137
$ws.Range("E11").Value = 0.9552813142152453
$ws.Range("F11").Value = 0.8941263331947911
$ws.Range("G11").Value = 0.9629231315451625
$ws.Range("H11").Value = 5.053231210095444
$ws.Range("I11").Value = 2.887776400614325

$ws.Range("A12").Value = "model_1_8_21"
$ws.Range("B12").Value = 0.9926796641800427
$ws.Range("C12").Value = 0.982910328232657
$ws.Range("D12").Value = 0.9213792753619734
$ws.Range("E12").Value = 0.9560664863804997
$ws.Range("F12").Value = 0.8828346793749129
$ws.Range("G12").Value = 0.957998606052838
$ws.Range("H12").Value = 4.951026291270453
$ws.Range("I12").Value = 2.83707272698324

$ws.Range("A13").Value = "model_1_8_20"
$ws.Range("B13").Value = 0.9927778109538639
$ws.Range("C13").Value = 0.9830053214387923
$ws.Range("D13").Value = 0.9231624992917499
$ws.Range("E13").Value = 0.9569285103000963
$ws.Range("F13").Value = 0.8709981492307853
$ws.Range("G13").Value = 0.9526735559113857
$ws.Range("H13").Value = 4.838730346401027
$ws.Range("I13").Value = 2.781406235714739

$ws.Range("A14").Value = "model_1_8_19"
$ws.Range("B14").Value = 0.9928796193748943
$ws.Range("C14").Value = 0.9831079126045191
$ws.Range("D14").Value = 0.9251193357064968
$ws.Range("E14").Value = 0.9578736594071373
$ws.Range("F14").Value = 0.8587200233430229
$ws.Range("G14").Value = 0.9469225856705474
$ws.Range("H14").Value = 4.715501406681466
$ws.Range("I14").Value = 2.720371810429702

$ws.Range("A15").Value = "model_1_8_7"
$ws.Range("B15").Value = 0.9929048023629121
$ws.Range("C15").Value = 0.984641808741168
$ws.Range("D15").Value = 0.9661156443635583
$ws.Range("E15").Value = 0.9773920908266507
$ws.Range("F15").Value = 0.8556829474903415
$ws.Range("G15").Value = 0.8609367118196798
$ws.Range("H15").Value = 2.133818231657955
$ws.Range("I15").Value = 1.459939741795528

$ws.Range("A16").Value = "model_1_8_18"
$ws.Range("B16").Value = 0.9929837865534251
$ws.Range("C16").Value = 0.9832182672592527
$ws.Range("D16").Value = 0.9272630761351848
$ws.Range("E16").Value = 0.9589081203003738
$ws.Range("F16").Value = 0.8461574305984528
$ws.Range("G16").Value = 0.9407364162200518
$ws.Range("H16").Value = 4.580502457320982
$ws.Range("I16").Value = 2.653569942207871

$ws.Range("A17").Value = "model_1_8_17"
$ws.Range("B17").Value = 0.9930885060802673
$ws.Range("C17").Value = 0.983336176764944
$ws.Range("D17").Value = 0.9296075527219122
$ws.Range("E17").Value = 0.9600382418850384
$ws.Range("F17").Value = 0.8335282244830754
$ws.Range("G17").Value = 0.9341267432181041
$ws.Range("H17").Value = 4.432862439073919
$ws.Range("I17").Value = 2.580590640943792

$ws.Range("A18").Value = "model_1_8_16"
$ws.Range("B18").Value = 0.9931911680882776
$ws.Range("C18").Value = 0.9834615521991642
$ws.Range("D18").Value = 0.9321656412728496
$ws.Range("E18").Value = 0.96126979076986
$ws.Range("F18").Value = 0.8211471557514141
$ws.Range("G18").Value = 0.9270985514042782
$ws.Range("H18").Value = 4.271770516690861
$ws.Range("I18").Value = 2.501061519204635

$ws.Range("A19").Value = "model_1_8_8"
$ws.Range("B19").Value = 0.9931969896930963
$ws.Range("C19").Value = 0.9845522136595691
$ws.Range("D19").Value = 0.9614088678938632
$ws.Range("E19").Value = 0.9751909434317682
$ws.Range("F19").Value = 0.8204450684770003
$ws.Range("G19").Value = 0.8659591583856167
$ws.Range("H19").Value = 2.43022066442468
$ws.Range("I19").Value = 1.602082145796641

$ws.Range("A20").Value = "model_1_8_15"
$ws.Range("B20").Value = 0.9932884519424908
$ws.Range("C20").Value = 0.9835938107097068
$ws.Range("D20").Value = 0.934951345678081
$ws.Range("E20").Value = 0.9626089275340313
$ws.Range("F20").Value = 0.8094146939689127
$ws.Range("G20").Value = 0.9196845138227902
$ws.Range("H20").Value = 4.09634481547139
$ws.Range("I20").Value = 2.414584748322265

$ws.Range("A21").Value = "model_1_8_14"
$ws.Range("B21").Value = 0.9933756765800372
$ws.Range("C21").Value = 0.9837320396781639
$ws.Range("D21").Value = 0.9379749485277848
$ws.Range("E21").Value = 0.964060000781202
$ws.Range("F21").Value = 0.7988953767113822
$ws.Range("G21").Value = 0.9119357892772911
$ws.Range("H21").Value = 3.90593780418822
$ws.Range("I21").Value = 2.320879510674811

$ws.Range("A22").Value = "model_1_8_9"
$ws.Range("B22").Value = 0.993378080348888
$ws.Range("C22").Value = 0.9844372628405375
$ws.Range("D22").Value = 0.9568941497030916
$ws.Range("E22").Value = 0.9730662950820839
$ws.Range("F22").Value = 0.7986054814723024
$ws.Range("G22").Value = 0.8724029757916056
$ws.Range("H22").Value = 2.714528505176578
$ws.Range("I22").Value = 1.73928450888385

$ws.Range("A23").Value = "model_1_8_13"
$ws.Range("B23").Value = 0.9934469481818605
$ws.Range("C23").Value = 0.983874798189598
$ws.Range("D23").Value = 0.9412473877712207
$ws.Range("E23").Value = 0.9656274009470578
$ws.Range("F23").Value = 0.7903000003117472
$ws.Range("G23").Value = 0.9039331513789232
$ws.Range("H23").Value = 3.699860681324865
$ws.Range("I23").Value = 2.219662287273748

$ws.Range("A24").Value = "model_1_8_10"
$ws.Range("B24").Value = 0.9934747406978288
$ws.Range("C24").Value = 0.9843060128769729
$ws.Range("D24").Value = 0.9526052515034771
$ws.Range("E24").Value = 0.9710377845625433
$ws.Range("F24").Value = 0.7869482146112975
$ws.Range("G24").Value = 0.879760477085427
$ws.Range("H24").Value = 2.984615658972711
$ws.Range("I24").Value = 1.870278626978532

$ws.Range("A25").Value = "model_1_8_12"
$ws.Range("B25").Value = 0.9934944123989866
$ws.Range("C25").Value = 0.9840201414283822
$ws.Range("D25").Value = 0.9447756116629622
$ws.Range("E25").Value = 0.9673132598534521
$ws.Range("F25").Value = 0.7845758015947896
$ws.Range("G25").Value = 0.895785620984536
$ws.Range("H25").Value = 3.477675890610645
$ws.Range("I25").Value = 2.110795412516202

$ws.Range("A26").Value = "model_1_8_11"
$ws.Range("B26").Value = 0.9935080006926251
$ws.Range("C26").Value = 0.9841651717109168
$ws.Range("D26").Value = 0.9485623713856401
$ws.Range("E26").Value = 0.9691177199597444
$ws.Range("F26").Value = 0.7829370493363405
$ws.Range("G26").Value = 0.8876556340313051
$ws.Range("H26").Value = 3.239210180303076
$ws.Range("I26").Value = 1.994269686874748
